$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$medium = -4138
$gray = 8355711
$orange = 32250
$greyFill = 15921906

function Set-Edge($range, $edgeIndex, $style, $weight, $color) {
    $range.Borders.Item($edgeIndex).LineStyle = $style
    if ($style -ne 0) {
        $range.Borders.Item($edgeIndex).Weight = $weight
        $range.Borders.Item($edgeIndex).Color = $color
    }
}

# Font/Fill/Alignment across the whole highlighted range
$full = $ws.Range("B3:F3")
$full.Font.Bold = $true
$full.Font.Color = $orange
$full.Interior.Pattern = 1
$full.Interior.Color = $greyFill
$full.VerticalAlignment = -4108
Write-Host "font/fill done"

# Borders per cell (must be done individually)
foreach ($addr in @("B3","C3","D3","E3","F3")) {
    $r = $ws.Range($addr)
    if ($addr -eq "B3") {
        Set-Edge $r 7 1 $medium $gray
    } else {
        Set-Edge $r 7 0 $medium $gray
    }
    Set-Edge $r 8 1 $medium $gray
    Set-Edge $r 9 1 $medium $gray
    Set-Edge $r 10 1 $medium $gray
}
Write-Host "borders done"

$ws.Range("F3").HorizontalAlignment = -4152
Write-Host "done"
